$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2, D2 updated; C2, E2 cleared out entirely)
$ws.Range("B2").Value = 10.582135398461524
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 14.565089373690022
$ws.Range("E2").ClearContents()

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 9.564422013072253
$ws.Range("C3").Value = -3.7209366071500796
$ws.Range("D3").Value = 14.367451801793614
$ws.Range("E3").Value = -3.2448504095349477

# Update the sheet's active selection to match the new data extent
$ws.Range("B1:E3").Select()
